$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (StreetTRACKS Gold Shares / GLD) - metric refresh, identity unchanged
$ws.Range("D2").Value = 387.17
$ws.Range("E2").Value = 70.59999999999999
$ws.Range("F2").Value = -0.18
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 83
$ws.Range("K2").Value = 60.4
$ws.Range("N2").Value = 50.60178744571824

# Row 3 - now identifies "Gold Feb 26" / "GC=F" (was Newmont Corporation / NEM)
$ws.Range("B3").Value = "Gold Feb 26"
$ws.Range("C3").Value = "GC=F"
$ws.Range("D3").Value = 4270.1
$ws.Range("E3").Value = 73.8
$ws.Range("F3").Value = 1.23
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 73
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 59.4
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 50.60178744571824

# Row 4 - now identifies "Newmont Corporation" / "NEM" (was Gold Feb 26 / GC=F)
$ws.Range("B4").Value = "Newmont Corporation"
$ws.Range("C4").Value = "NEM"
$ws.Range("D4").Value = 90.42
$ws.Range("E4").Value = 56.6
$ws.Range("F4").Value = -0.34
$ws.Range("H4").Value = 76
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 83
$ws.Range("K4").Value = 58.2
$ws.Range("N4").Value = 50.60178744571824
